$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("A44").Value = "Stephan Morbitzer"
$ws.Range("B44").Value = "Zahide Jashari (Optional Morbitzer)"

# Restore the view state (scroll position / active selection) that Excel
# records after the edit.
$win = $excel.ActiveWindow
$win.ScrollRow = 15
$win.ScrollColumn = 1
$ws.Range("B45").Select()
